$d = $word.ActiveDocument

$replacements = @(
    @("65×14=", "45×93="),
    @("75×44=", "39×79="),
    @("78×77=", "34×80="),
    @("73×48=", "27×85="),
    @("63×90=", "96×66="),
    @("93×54=", "45×23="),
    @("89×33=", "79×89="),
    @("60×87=", "60×98="),
    @("42×53=", "79×90="),
    @("15×75=", "61×45="),
    @("43×22=", "55×35="),
    @("61×73=", "51×68="),
    @("64×13=", "88×24="),
    @("54×95=", "79×16="),
    @("51×73=", "33×18="),
    @("56×22=", "66×63="),
    @("90×21=", "12×70="),
    @("84×21=", "29×87="),
    @("13×28=", "77×56="),
    @("74×55=", "25×15="),
    @("29×71=", "17×60="),
    @("98×57=", "41×84="),
    @("21×31=", "58×71="),
    @("76×68=", "72×32="),
    @("85×88=", "78×84=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
